$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume figures.
# Values are set with a leading apostrophe to force plain text (preventing Excel
# from re-interpreting numeric-looking / percent-looking strings as numbers),
# then ClearFormats() removes the resulting quote-prefix style so the cell keeps
# the workbook default (unstyled) formatting, matching the original cells.

$ws.Range("D2").Value = "'26.341.39"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +1.24%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.683.49"
$ws.Range("D3").ClearFormats()
$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  +0.35%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'218.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +0.65%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'0.5530"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'1.008"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +0.32%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.2701"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +1.69%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.06503"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +1.50%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'22.11"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +1.34%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.07542"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +1.29%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'4.549"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +0.85%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'1.679.11"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +0.63%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.5815"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -0.47%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.000008459"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -1.41%  "
$ws.Range("E15").ClearFormats()
$ws.Range("E16").Value = "'  +1.18%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'26.380.20"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +1.08%  "
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = "'  -0.04%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  +0.33%  "
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = "'  +1.27%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'191.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -0.38%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'6.233"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +0.38%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D24").Value = "'147.42"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +1.82%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.1330"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +11.20%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'7.902"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +3.66%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'15.82"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +0.85%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'0.06350"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -2.33%  "
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = "'  +3.99%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  +0.48%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'  +1.52%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  +1.83%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'1.670"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +1.23%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'1.040"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +2.02%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'0.6217"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +1.60%  "
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = "'  +1.32%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'2.722"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +1.48%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'6.235"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -0.56%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'1.111.97"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +2.03%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D41").Value = "'0.8731"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +1.10%  "
$ws.Range("E41").ClearFormats()
$ws.Range("E42").Value = "'  +0.68%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'100.79"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  -0.13%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'1.832.62"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +0.93%  "
$ws.Range("E44").ClearFormats()
$ws.Range("E45").Value = "'  -2.34%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'57.35"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +1.52%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'8.207"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +1.57%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'1.004"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -0.38%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.05274"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +0.79%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.4295"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +0.22%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'6.080"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +0.17%  "
$ws.Range("E51").ClearFormats()
